# Update countries & provincias Spain
# - Swap the (name, data) pairing for three pairs of adjacent rows so the
#   country names appear in their new order, and refresh several rows of
#   case-count figures with newer data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (Brasil): refreshed totals ---------------------------------
$ws.Range("B14").Value = 29015
$ws.Range("C14").Value = 405
$ws.Range("E14").Value = 13229

# --- Rows 46/47: Catar now listed before Panama -------------------------
# Row 46 keeps its position but becomes "Catar" with brand-new figures;
# row 47 becomes "Panama" and inherits the figures that used to sit in
# row 46 (i.e. the data shifts down one row together with the name swap).
$ws.Range("A46").Value = "Catar"
$ws.Range("B46").Value = 4103
$ws.Range("C46").Value = 392
$ws.Range("D46").Value = 415
$ws.Range("E46").Value = 3681
$ws.Range("F46").Value = 37
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 7

$ws.Range("A47").Value = "Panama"
$ws.Range("B47").Value = 3751
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 75
$ws.Range("E47").Value = 3573
$ws.Range("F47").Value = 106
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 103

# --- Rows 70/71: Kazajistan now listed before Uzbekistan -----------------
$ws.Range("A70").Value = "Kazajistan"
$ws.Range("B70").Value = 1362
$ws.Range("C70").Value = 67
$ws.Range("D70").Value = 273
$ws.Range("E70").Value = 1072
$ws.Range("F70").Value = 22
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 17

$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("B71").Value = 1349
$ws.Range("C71").Value = 47
$ws.Range("D71").Value = 107
$ws.Range("E71").Value = 1238
$ws.Range("F71").Value = 8
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 4

# --- Rows 74/75: Bosnia y Herzegovina now listed before Armenia ----------
$ws.Range("A74").Value = "Bosnia y Herzegovina"
$ws.Range("B74").Value = 1167
$ws.Range("C74").Value = 57
$ws.Range("D74").Value = 277
$ws.Range("E74").Value = 847
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 43

$ws.Range("A75").Value = "Armenia"
$ws.Range("B75").Value = 1159
$ws.Range("C75").Value = 48
$ws.Range("D75").Value = 358
$ws.Range("E75").Value = 783
$ws.Range("F75").Value = 30
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 18

# --- Row 79 (Eslovaquia): refreshed totals -------------------------------
$ws.Range("D79").Value = 167
$ws.Range("E79").Value = 802
